$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B (Coin name) updates
$ws.Cells.Item(23, 2).Value = "InternetComputer(DFINITY)"
$ws.Cells.Item(24, 2).Value = "Litecoin"
$ws.Cells.Item(42, 2).Value = "WEMIXToken"
$ws.Cells.Item(43, 2).Value = "VeChain"
$ws.Cells.Item(49, 2).Value = "BabyDogeCoin"
$ws.Cells.Item(50, 2).Value = "Monero"
$ws.Cells.Item(51, 2).Value = "Maker"

# Column C (Link) updates
$ws.Cells.Item(23, 3).Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Cells.Item(24, 3).Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"

# Column D (Price) updates - force text to avoid numeric auto-conversion,
# then clear the format so no extra style index is introduced.
$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = "69.336.32"
$ws.Cells.Item(2, 4).ClearFormats()
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = "3.897.63"
$ws.Cells.Item(3, 4).ClearFormats()
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "528.92"
$ws.Cells.Item(5, 4).ClearFormats()
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "144.51"
$ws.Cells.Item(6, 4).ClearFormats()
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "0.611"
$ws.Cells.Item(7, 4).ClearFormats()
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "0.0000336"
$ws.Cells.Item(11, 4).ClearFormats()
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "42.05"
$ws.Cells.Item(12, 4).ClearFormats()
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "4.518.78"
$ws.Cells.Item(13, 4).ClearFormats()
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "10.24"
$ws.Cells.Item(14, 4).ClearFormats()
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "3.897.43"
$ws.Cells.Item(15, 4).ClearFormats()
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "19.80"
$ws.Cells.Item(19, 4).ClearFormats()
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "69.218.54"
$ws.Cells.Item(20, 4).ClearFormats()
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "425.94"
$ws.Cells.Item(21, 4).ClearFormats()
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "3.39"
$ws.Cells.Item(22, 4).ClearFormats()
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "14.13"
$ws.Cells.Item(23, 4).ClearFormats()
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "88.12"
$ws.Cells.Item(24, 4).ClearFormats()
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "11.37"
$ws.Cells.Item(26, 4).ClearFormats()
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "10.58"
$ws.Cells.Item(27, 4).ClearFormats()
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "36.36"
$ws.Cells.Item(28, 4).ClearFormats()
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "689.92"
$ws.Cells.Item(29, 4).ClearFormats()
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "13.19"
$ws.Cells.Item(30, 4).ClearFormats()
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "68.81"
$ws.Cells.Item(33, 4).ClearFormats()
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "0.0₃0889"
$ws.Cells.Item(34, 4).ClearFormats()
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "3.30"
$ws.Cells.Item(40, 4).ClearFormats()
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "3.22"
$ws.Cells.Item(42, 4).ClearFormats()
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "0.0480"
$ws.Cells.Item(43, 4).ClearFormats()
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "2.81"
$ws.Cells.Item(44, 4).ClearFormats()
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "3.39"
$ws.Cells.Item(45, 4).ClearFormats()
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "0.000290"
$ws.Cells.Item(46, 4).ClearFormats()
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "0.0₆0347"
$ws.Cells.Item(49, 4).ClearFormats()
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "146.19"
$ws.Cells.Item(50, 4).ClearFormats()
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "2.743.25"
$ws.Cells.Item(51, 4).ClearFormats()

# Column E (Volume/percent change) updates
$ws.Cells.Item(2, 5).Value = "  +1.66%  "
$ws.Cells.Item(3, 5).Value = "  +0.26%  "
$ws.Cells.Item(4, 5).Value = "  -0.03%  "
$ws.Cells.Item(5, 5).Value = "  +9.60%  "
$ws.Cells.Item(6, 5).Value = "  -0.75%  "
$ws.Cells.Item(7, 5).Value = "  -1.73%  "
$ws.Cells.Item(8, 5).Value = "  +0.02%  "
$ws.Cells.Item(9, 5).Value = "  -3.54%  "
$ws.Cells.Item(10, 5).Value = "  -4.84%  "
$ws.Cells.Item(11, 5).Value = "  -5.56%  "
$ws.Cells.Item(12, 5).Value = "  -2.43%  "
$ws.Cells.Item(13, 5).Value = "  +0.29%  "
$ws.Cells.Item(14, 5).Value = "  -2.67%  "
$ws.Cells.Item(15, 5).Value = "  +0.56%  "
$ws.Cells.Item(16, 5).Value = "  -1.97%  "
$ws.Cells.Item(17, 5).Value = "  -1.25%  "
$ws.Cells.Item(18, 5).Value = "  +6.66%  "
$ws.Cells.Item(19, 5).Value = "  -0.79%  "
$ws.Cells.Item(20, 5).Value = "  +1.41%  "
$ws.Cells.Item(22, 5).Value = "  -5.67%  "
$ws.Cells.Item(23, 5).Value = "  -4.60%  "
$ws.Cells.Item(24, 5).Value = "  -0.73%  "
$ws.Cells.Item(25, 5).Value = "  +9.89%  "
$ws.Cells.Item(26, 5).Value = "  -9.71%  "
$ws.Cells.Item(27, 5).Value = "  -3.94%  "
$ws.Cells.Item(28, 5).Value = "  -2.44%  "
$ws.Cells.Item(29, 5).Value = "  -4.14%  "
$ws.Cells.Item(30, 5).Value = "  -2.29%  "
$ws.Cells.Item(31, 5).Value = "  -3.27%  "
$ws.Cells.Item(32, 5).Value = "  -2.91%  "
$ws.Cells.Item(33, 5).Value = "  +11.45%  "
$ws.Cells.Item(34, 5).Value = "  +1.63%  "
$ws.Cells.Item(35, 5).Value = "  +8.14%  "
$ws.Cells.Item(36, 5).Value = "  -2.14%  "
$ws.Cells.Item(38, 5).Value = "  +1.93%  "
$ws.Cells.Item(39, 5).Value = "  +0.02%  "
$ws.Cells.Item(40, 5).Value = "  +7.45%  "
$ws.Cells.Item(41, 5).Value = "  -0.11%  "
$ws.Cells.Item(42, 5).Value = "  +8.97%  "
$ws.Cells.Item(43, 5).Value = "  -3.63%  "
$ws.Cells.Item(44, 5).Value = "  -6.70%  "
$ws.Cells.Item(45, 5).Value = "  +0.74%  "
$ws.Cells.Item(46, 5).Value = "  +19.56%  "
$ws.Cells.Item(47, 5).Value = "  -1.69%  "
$ws.Cells.Item(48, 5).Value = "  +6.75%  "
$ws.Cells.Item(49, 5).Value = "  -4.13%  "
$ws.Cells.Item(50, 5).Value = "  +1.14%  "
$ws.Cells.Item(51, 5).Value = "  +14.77%  "
